$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix existing rows 66 and 67 (dates shift back by one day)
$ws.Range("E66").Value = 43839
$ws.Range("E67").Value = 43840

# Add new row 68
$ws.Range("E68").Value = 43841
$ws.Range("F68").Value = 5
$ws.Range("G68").Value = "Darstellung von allgemeinen Informationen über stream im Menüpunkt 'About'"

# Add new row 69
$ws.Range("E69").Value = 43842
$ws.Range("F69").Value = 1
$ws.Range("G69").Value = "Implementieren 2 neuer Sprachen (Französisch und Spanisch) durch ResourceBundles"

# Copy formatting of row 67 (E:G) down into the newly added rows 68 and 69
$ws.Range("E67:G67").Copy()
$ws.Range("E68:G68").PasteSpecial(-4122)
$ws.Range("E69:G69").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update view - top-left cell and selection
$win = $excel.ActiveWindow
$win.SetTopLeftVisibleCell($ws.Range("B64"))
$ws.Range("E70").Select()

# Recalculate formulas (C5 = SUM(F:F) must reflect the new rows)
$excel.Calculate()

$wb.Save()
